$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark a couple more tasks as "done"
$ws.Range("B2").Value = "done"
$ws.Range("B4").Value = "done"

# New "routing directions" screen task got a second VC (carousel) — the
# login task is now superseded, so strike it through…
$ws.Range("B9").Font.Strikethrough = $true

# …and mark the two VC items below it ("user info" and "directions -
# carousel") as done in a new column C.
$ws.Range("C10").Value = "done"
$ws.Range("C11").Value = "done"

# Give the new column B text room (matches the wider title column once the
# longer VC strings live there).
$ws.Columns.Item(2).ColumnWidth = 68.5

# Leave the selection where the edit finished.
$ws.Range("C14").Select() | Out-Null
